$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2135922330097087
$ws.Range("C2").Value = 0.5145631067961165
$ws.Range("J2").Value = 0.01456310679611651
$ws.Range("P2").Value = 0.1504854368932039
$ws.Range("S2").Value = 0.1067961165048544
$ws.Range("B3").Value = 0.01801801801801802
$ws.Range("C3").Value = 0.05405405405405406
$ws.Range("J3").Value = 0.01801801801801802
$ws.Range("P3").Value = 0.7477477477477478
$ws.Range("S3").Value = 0.1621621621621622
$ws.Range("O4").Value = 0.04347826086956522
$ws.Range("P4").Value = 0.7826086956521739
$ws.Range("S4").Value = 0.1739130434782609
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.07109004739336493
$ws.Range("D6").Value = 0.004739336492890996
$ws.Range("F6").Value = 0.08056872037914692
$ws.Range("J6").Value = 0.2464454976303317
$ws.Range("O6").Value = 0.01421800947867299
$ws.Range("Q6").Value = 0.1137440758293839
$ws.Range("R6").Value = 0.07109004739336493
$ws.Range("S6").Value = 0.3981042654028436
$ws.Range("B7").Value = 0.08843537414965986
$ws.Range("D7").Value = 0.02040816326530612
$ws.Range("F7").Value = 0.08843537414965986
$ws.Range("J7").Value = 0.06802721088435375
$ws.Range("Q7").Value = 0.163265306122449
$ws.Range("R7").Value = 0.08163265306122448
$ws.Range("S7").Value = 0.4897959183673469
$ws.Range("B8").Value = 0.06542056074766354
$ws.Range("D8").Value = 0.006230529595015576
$ws.Range("F8").Value = 0.06230529595015576
$ws.Range("J8").Value = 0.1090342679127726
$ws.Range("O8").Value = 0.01246105919003115
$ws.Range("Q8").Value = 0.1713395638629283
$ws.Range("R8").Value = 0.1090342679127726
$ws.Range("S8").Value = 0.4641744548286604
$ws.Range("B9").Value = 0.0660377358490566
$ws.Range("D9").Value = 0.02358490566037736
$ws.Range("F9").Value = 0.0660377358490566
$ws.Range("J9").Value = 0.07075471698113207
$ws.Range("O9").Value = 0.009433962264150943
$ws.Range("Q9").Value = 0.169811320754717
$ws.Range("R9").Value = 0.09433962264150944
$ws.Range("S9").Value = 0.5
$ws.Range("B10").Value = 0.0954356846473029
$ws.Range("D10").Value = 0.01244813278008299
$ws.Range("E10").Value = 0.002074688796680498
$ws.Range("F10").Value = 0.09336099585062241
$ws.Range("J10").Value = 0.1026970954356846
$ws.Range("O10").Value = 0.01452282157676349
$ws.Range("Q10").Value = 0.1970954356846473
$ws.Range("R10").Value = 0.09958506224066389
$ws.Range("S10").Value = 0.3827800829875518
$ws.Range("G11").Value = 0.1122994652406417
$ws.Range("J11").Value = 0.05882352941176471
$ws.Range("K11").Value = 0.1336898395721925
$ws.Range("L11").Value = 0.6898395721925134
$ws.Range("S11").Value = 0.0053475935828877
$ws.Range("G12").Value = 0.781021897810219
$ws.Range("J12").Value = 0.1167883211678832
$ws.Range("K12").Value = 0.0291970802919708
$ws.Range("L12").Value = 0.0583941605839416
$ws.Range("S12").Value = 0.0145985401459854
$ws.Range("G13").Value = 0.7777777777777778
$ws.Range("J13").Value = 0.1851851851851852
$ws.Range("S13").Value = 0.03703703703703703
$ws.Range("F15").Value = 0.01298701298701299
$ws.Range("H15").Value = 0.1233766233766234
$ws.Range("I15").Value = 0.1168831168831169
$ws.Range("J15").Value = 0.4090909090909091
$ws.Range("K15").Value = 0.07792207792207792
$ws.Range("M15").Value = 0.01948051948051948
$ws.Range("O15").Value = 0.07142857142857142
$ws.Range("S15").Value = 0.1688311688311688
$ws.Range("F16").Value = 0.01550387596899225
$ws.Range("H16").Value = 0.1395348837209302
$ws.Range("I16").Value = 0.09302325581395349
$ws.Range("J16").Value = 0.4418604651162791
$ws.Range("K16").Value = 0.1007751937984496
$ws.Range("M16").Value = 0.01550387596899225
$ws.Range("N16").Value = 0.007751937984496124
$ws.Range("O16").Value = 0.03875968992248062
$ws.Range("S16").Value = 0.1472868217054264
$ws.Range("F17").Value = 0.00911854103343465
$ws.Range("H17").Value = 0.1914893617021277
$ws.Range("I17").Value = 0.1124620060790274
$ws.Range("J17").Value = 0.4437689969604863
$ws.Range("K17").Value = 0.05167173252279635
$ws.Range("M17").Value = 0.02127659574468085
$ws.Range("N17").Value = 0.00303951367781155
$ws.Range("O17").Value = 0.060790273556231
$ws.Range("S17").Value = 0.1063829787234043
$ws.Range("F18").Value = 0.01685393258426966
$ws.Range("H18").Value = 0.1853932584269663
$ws.Range("I18").Value = 0.1123595505617977
$ws.Range("J18").Value = 0.3876404494382023
$ws.Range("K18").Value = 0.07865168539325842
$ws.Range("M18").Value = 0.01685393258426966
$ws.Range("N18").Value = 0.005617977528089887
$ws.Range("O18").Value = 0.05056179775280899
$ws.Range("S18").Value = 0.1460674157303371
$ws.Range("F19").Value = 0.01653696498054475
$ws.Range("H19").Value = 0.183852140077821
$ws.Range("I19").Value = 0.1206225680933852
$ws.Range("J19").Value = 0.3920233463035019
$ws.Range("K19").Value = 0.0963035019455253
$ws.Range("M19").Value = 0.01361867704280156
$ws.Range("N19").Value = 0.0009727626459143969
$ws.Range("O19").Value = 0.07101167315175097
$ws.Range("S19").Value = 0.1050583657587549
